# Excel COM-interop script: "added additional sentiment analysis"
# - Recalculates the back-adjusted "Adj Close" (column F) values for the
#   existing historical rows (the adjustment factor shifts slightly whenever
#   new trading days are appended to the dataset).
# - Appends two new rows (149-150) of OHLCV + Adj Close data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "Adj Close" column (F) for rows whose back-adjustment factor
#    changed as a result of the newly appended price history below.
# ---------------------------------------------------------------------------
$adjCloseUpdates = @{
  4 = 79.81497192382812
  5 = 75.43687438964844
  6 = 74.96308135986328
  7 = 74.647216796875
  8 = 72.51519012451172
  11 = 72.12037658691406
  12 = 67.33868408203125
  14 = 63.1268424987793
  15 = 63.45415115356445
  16 = 67.45265960693359
  17 = 72.97271728515625
  18 = 73.830810546875
  19 = 72.03501892089844
  21 = 67.58535003662109
  23 = 67.55881500244141
  25 = 72.76924896240234
  26 = 67.49688720703125
  29 = 67.04367065429688
  31 = 68.41683197021484
  32 = 71.21662139892578
  33 = 68.74673461914062
  34 = 71.54653167724609
  35 = 72.24203491210938
  38 = 76.37040710449219
  40 = 75.32717132568359
  41 = 76.51306915283203
  42 = 74.48291778564453
  45 = 71.69865417480469
  46 = 71.69865417480469
  47 = 67.54925537109375
  50 = 72.25552368164062
  52 = 73.59375
  54 = 71.42922210693359
  56 = 74.04818725585938
  57 = 75.16162872314453
  58 = 75.47846221923828
  60 = 69.76641845703125
  61 = 68.61677551269531
  62 = 67.892578125
  65 = 66.19979858398438
  66 = 67.53050231933594
  68 = 70.78119659423828
  70 = 70.39765930175781
  74 = 68.03251647949219
  75 = 71.60307312011719
  76 = 71.49348449707031
  77 = 72.30623626708984
  79 = 70.92481994628906
  81 = 72.99199676513672
  83 = 74.51708984375
  84 = 75.86843109130859
  85 = 75.27962493896484
  88 = 74.19856262207031
  89 = 72.81824493408203
  90 = 74.05376434326172
  92 = 73.64837646484375
  93 = 74.57499694824219
  94 = 70.89739990234375
  98 = 79.2752685546875
  99 = 78.52536010742188
  102 = 85.75168609619141
  105 = 78.59353637695312
  108 = 70.72444915771484
  111 = 74.45298004150391
  120 = 75.03302001953125
  122 = 76.94031524658203
}

foreach ($row in $adjCloseUpdates.Keys) {
    $ws.Cells.Item($row, 6).Value = $adjCloseUpdates[$row]
}

# ---------------------------------------------------------------------------
# 2) Append two new rows of data (149 and 150), matching the style/format of
#    the last existing data row (148).
# ---------------------------------------------------------------------------
$ws.Range("A148:G148").Copy()
$ws.Range("A149:G150").PasteSpecial(-4122)

$newRows = New-Object 'object[,]' 2,7
$newRows[0,0] = 44809
$newRows[0,1] = 86.58000183105469
$newRows[0,2] = 88.16999816894531
$newRows[0,3] = 85.06999969482422
$newRows[0,4] = 87.33999633789062
$newRows[0,5] = 87.33999633789062
$newRows[0,6] = 29955300

$newRows[1,0] = 44813
$newRows[1,1] = 87.58999633789062
$newRows[1,2] = 88.15499877929688
$newRows[1,3] = 87.26999664306641
$newRows[1,4] = 87.33999633789062
$newRows[1,5] = 87.33999633789062
$newRows[1,6] = 7481322

$ws.Range("A149:G150").Value = $newRows
